$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Couples")
$ws.Range("B2").Value = 0.81926267736065861
$ws.Range("B3").Value = 0.81816541766195372
$ws.Range("B4").Value = 0.81263969347138987
$ws.Range("B5").Value = 0.82765408852075384
$ws.Range("B6").Value = 0.82950872081543936
$ws.Range("B7").Value = 0.83611248688842976
$ws.Range("B8").Value = 0.83775723517092482
$ws.Range("B9").Value = 0.84284125820157529
$ws.Range("B10").Value = 0.84858463562677389
$ws.Range("B11").Value = 0.86186749844364374
$ws.Range("B12").Value = 0.8771080374304252
$ws.Range("B13").Value = 0.87584419058033092
$ws.Range("B14").Value = 0.88957720844980337

$ws = $wb.Worksheets.Item("SingleAC_Females")
$ws.Range("B2").Value = 0.3991450072596518
$ws.Range("B3").Value = 0.41948999821695471
$ws.Range("B4").Value = 0.39157863662299497
$ws.Range("B5").Value = 0.41433692743405076
$ws.Range("B6").Value = 0.41124963262829567
$ws.Range("B7").Value = 0.42607921647394814
$ws.Range("B8").Value = 0.45788102982282675
$ws.Range("B9").Value = 0.43475036545320306
$ws.Range("B10").Value = 0.4352103131468209
$ws.Range("B11").Value = 0.48020557694901034
$ws.Range("B12").Value = 0.45212571734411705
$ws.Range("B13").Value = 0.43647713649597536
$ws.Range("B14").Value = 0.43758876928331342

$ws = $wb.Worksheets.Item("SingleAC_Males")
$ws.Range("B2").Value = 0.53742812155347619
$ws.Range("B3").Value = 0.53844044474943309
$ws.Range("B4").Value = 0.52988290374081737
$ws.Range("B5").Value = 0.54149128104056998
$ws.Range("B6").Value = 0.55716858320426854
$ws.Range("B7").Value = 0.6085817810490497
$ws.Range("B8").Value = 0.61658657335201739
$ws.Range("B9").Value = 0.59176278878693189
$ws.Range("B10").Value = 0.58548745841077365
$ws.Range("B11").Value = 0.58854728020588165
$ws.Range("B12").Value = 0.62663450587764535
$ws.Range("B13").Value = 0.61267069061990309
$ws.Range("B14").Value = 0.6023965824847306

$ws = $wb.Worksheets.Item("SingleDep_Females")
$ws.Range("B2").Value = 0.30960911864445839
$ws.Range("B3").Value = 0.3009604726704741
$ws.Range("B4").Value = 0.30537970485383314
$ws.Range("B5").Value = 0.27890116565669643
$ws.Range("B6").Value = 0.27448488810723043
$ws.Range("B7").Value = 0.30109497280240638
$ws.Range("B8").Value = 0.30056498819528676
$ws.Range("B9").Value = 0.30681190730380958
$ws.Range("B10").Value = 0.35913699591968529
$ws.Range("B11").Value = 0.34448547398830465
$ws.Range("B12").Value = 0.35897695098347643
$ws.Range("B13").Value = 0.3588495010998119
$ws.Range("B14").Value = 0.36154285650172741

$ws = $wb.Worksheets.Item("SingleDep_Males")
$ws.Range("B2").Value = 0.41638526322582725
$ws.Range("B3").Value = 0.4000992497086247
$ws.Range("B4").Value = 0.39148649281255121
$ws.Range("B5").Value = 0.39413842587482062
$ws.Range("B6").Value = 0.40679199236816832
$ws.Range("B7").Value = 0.40240254301445588
$ws.Range("B8").Value = 0.41093893006058185
$ws.Range("B9").Value = 0.42075425802135913
$ws.Range("B10").Value = 0.42573444595703996
$ws.Range("B11").Value = 0.40895011068374643
$ws.Range("B12").Value = 0.43199059291587072
$ws.Range("B13").Value = 0.43794253919446102
$ws.Range("B14").Value = 0.45836433725671061

$ws = $wb.Worksheets.Item("Single_female")
$ws.Range("B2").Value = 0.29890599073681112
$ws.Range("B3").Value = 0.30143439406655592
$ws.Range("B4").Value = 0.30152689696624108
$ws.Range("B5").Value = 0.31000832388467597
$ws.Range("B6").Value = 0.30977756771418391
$ws.Range("B7").Value = 0.31115539905180917
$ws.Range("B8").Value = 0.30826073273200044
$ws.Range("B9").Value = 0.30027367523619586
$ws.Range("B10").Value = 0.3114533855079103
$ws.Range("B11").Value = 0.31075908202123842
$ws.Range("B12").Value = 0.31838468917312085
$ws.Range("B13").Value = 0.34228404069903945
$ws.Range("B14").Value = 0.33890969837670926

$ws = $wb.Worksheets.Item("Single_male")
$ws.Range("B2").Value = 0.44718650879882615
$ws.Range("B3").Value = 0.44121448255095902
$ws.Range("B4").Value = 0.44471308403527815
$ws.Range("B5").Value = 0.4533156680533294
$ws.Range("B6").Value = 0.47460676585552064
$ws.Range("B7").Value = 0.46812326004996629
$ws.Range("B8").Value = 0.48024505747805785
$ws.Range("B9").Value = 0.48649978080337819
$ws.Range("B10").Value = 0.48701374677103965
$ws.Range("B11").Value = 0.49963984755210444
$ws.Range("B12").Value = 0.52688230019279025
$ws.Range("B13").Value = 0.53298706390936279
$ws.Range("B14").Value = 0.51857282738817789
